$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.593669337425428
$ws.Range("C2").Value = 6.459742127343147
$ws.Range("E2").Value = 16.53736111016295
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 20.20755339058882
$ws.Range("H2").Value = 12.30318198815509
$ws.Range("I2").Value = 17.37289335160082
$ws.Range("K2").Value = 9.239354112517319
$ws.Range("O2").Value = 17.55997275444416

$ws.Range("B3").Value = 8.170753776752328
$ws.Range("C3").Value = 6.285196716249302
$ws.Range("E3").Value = 15.59869541374222
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 20.38590023328348
$ws.Range("H3").Value = 12.36479282231346
$ws.Range("I3").Value = 17.49469821902562
$ws.Range("K3").Value = 8.923135850748967
$ws.Range("O3").Value = 17.67549190488816

$ws.Range("B4").Value = 7.899586619791833
$ws.Range("C4").Value = 6.174879816275154
$ws.Range("E4").Value = 14.99705907009783
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 20.50464095743429
$ws.Range("H4").Value = 12.40479087960728
$ws.Range("I4").Value = 17.57338720478943
$ws.Range("K4").Value = 8.721760339334384
$ws.Range("O4").Value = 17.75090962806147

$ws.Range("B5").Value = 7.78631206904909
$ws.Range("C5").Value = 6.129180095764659
$ws.Range("E5").Value = 14.74579530734558
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 20.55533321542663
$ws.Range("H5").Value = 12.42163634939939
$ws.Range("I5").Value = 17.60643617724406
$ws.Range("K5").Value = 8.637965723695967
$ws.Range("O5").Value = 17.78277036789328

$ws.Range("B6").Value = 7.767339382544004
$ws.Range("C6").Value = 6.12154813276558
$ws.Range("E6").Value = 14.70371407509736
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 20.56388929532758
$ws.Range("H6").Value = 12.4244665176895
$ws.Range("I6").Value = 17.61198333266621
$ws.Range("K6").Value = 8.623949628769564
$ws.Range("O6").Value = 17.78812887679989

$ws.Range("B7").Value = 7.898070010391931
$ws.Range("C7").Value = 6.174266445467326
$ws.Range("E7").Value = 14.99369471305475
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 20.50531530365419
$ws.Range("H7").Value = 12.40501585224835
$ws.Range("I7").Value = 17.57382893398889
$ws.Range("K7").Value = 8.720637155516709
$ws.Range("O7").Value = 17.75133474995154

$ws.Range("B8").Value = 8.450299174778579
$ws.Range("C8").Value = 6.400238346091272
$ws.Range("E8").Value = 16.21910351277491
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 20.26711987467334
$ws.Range("H8").Value = 12.32397562863803
$ws.Range("I8").Value = 17.41408324891855
$ws.Range("K8").Value = 9.131866342995053
$ws.Range("O8").Value = 17.59887139140611

$ws.Range("B9").Value = 9.441382406846031
$ws.Range("C9").Value = 6.816557071391833
$ws.Range("E9").Value = 18.51012740470379
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 19.87410293913234
$ws.Range("H9").Value = 12.18223729891341
$ws.Range("I9").Value = 17.13168643409617
$ws.Range("K9").Value = 9.877857249889027
$ws.Range("O9").Value = 17.33556367659569

$ws.Range("B10").Value = 10.16174458786111
$ws.Range("C10").Value = 7.103841505207825
$ws.Range("E10").Value = 20.15008337421779
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 19.63164089995034
$ws.Range("H10").Value = 12.08853877584501
$ws.Range("I10").Value = 16.9429078509335
$ws.Range("K10").Value = 10.3854215066119
$ws.Range("O10").Value = 17.16394116581485

$ws.Range("B11").Value = 10.47155207736119
$ws.Range("C11").Value = 7.230092529560977
$ws.Range("E11").Value = 20.85367963826348
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 19.53163697325141
$ws.Range("H11").Value = 12.04817162327284
$ws.Range("I11").Value = 16.86106356455551
$ws.Range("K11").Value = 10.60689873039724
$ws.Range("O11").Value = 17.0906258149521

$ws.Range("B12").Value = 10.58628301183519
$ws.Range("C12").Value = 7.277232893459064
$ws.Range("E12").Value = 21.11404780258093
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 19.49526918990421
$ws.Range("H12").Value = 12.03320970117831
$ws.Range("I12").Value = 16.83064941218857
$ws.Range("K12").Value = 10.68936842927858
$ws.Range("O12").Value = 17.06354921121915

$ws.Range("B13").Value = 10.56168886112699
$ws.Range("C13").Value = 7.267110538891091
$ws.Range("E13").Value = 21.05824223985627
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 19.50303451866958
$ws.Range("H13").Value = 12.03641760211409
$ws.Range("I13").Value = 16.83717393505312
$ws.Range("K13").Value = 10.67166997371169
$ws.Range("O13").Value = 17.06935007906567

$ws.Range("B14").Value = 10.48104296100208
$ws.Range("C14").Value = 7.233984349135699
$ws.Range("E14").Value = 20.8752217190855
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 19.52861475755424
$ws.Range("H14").Value = 12.04693420011453
$ws.Range("I14").Value = 16.85854978710775
$ws.Range("K14").Value = 10.61371183476831
$ws.Range("O14").Value = 17.08838443758162

$ws.Range("B15").Value = 10.43130784811779
$ws.Range("C15").Value = 7.213605741115385
$ws.Range("E15").Value = 20.76232703529229
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 19.5444795750596
$ws.Range("H15").Value = 12.05341813828407
$ws.Range("I15").Value = 16.87171841430189
$ws.Range("K15").Value = 10.5780273892894
$ws.Range("O15").Value = 17.10013298354689

$ws.Range("B16").Value = 10.1411366645837
$ws.Range("C16").Value = 7.095498782546337
$ws.Range("E16").Value = 20.10325073864177
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 19.6383853399034
$ws.Range("H16").Value = 12.09122225167254
$ws.Range("I16").Value = 16.9483375899577
$ws.Range("K16").Value = 10.37075399810438
$ws.Range("O16").Value = 17.16882840798296

$ws.Range("B17").Value = 9.958532790923433
$ws.Range("C17").Value = 7.0218853000718
$ws.Range("E17").Value = 19.68808524919545
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 19.6986449638859
$ws.Range("H17").Value = 12.1149916559879
$ws.Range("I17").Value = 16.99637273028906
$ws.Range("K17").Value = 10.24115291728024
$ws.Range("O17").Value = 17.21219072474369

$ws.Range("B18").Value = 9.85182059969009
$ws.Range("C18").Value = 6.979128834326342
$ws.Range("E18").Value = 19.44529758798969
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 19.73427226157355
$ws.Range("H18").Value = 12.12887562160819
$ws.Range("I18").Value = 17.02438086084724
$ws.Range("K18").Value = 10.16572578037976
$ws.Range("O18").Value = 17.23757914744944

$ws.Range("B19").Value = 9.815401085529267
$ws.Range("C19").Value = 6.96458176375305
$ws.Range("E19").Value = 19.36240663706594
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 19.74650058409296
$ws.Range("H19").Value = 12.13361298844103
$ws.Range("I19").Value = 17.03392916789521
$ws.Range("K19").Value = 10.14003708095493
$ws.Range("O19").Value = 17.24625203180155

$ws.Range("B20").Value = 9.978145613486337
$ws.Range("C20").Value = 7.02976485770832
$ws.Range("E20").Value = 19.73269360156692
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 19.69212993158721
$ws.Range("H20").Value = 12.11243937841996
$ws.Range("I20").Value = 16.99122003121271
$ws.Range("K20").Value = 10.25504101039869
$ws.Range("O20").Value = 17.20752839425113

$ws.Range("B21").Value = 10.504800901668
$ws.Range("C21").Value = 7.243732673848198
$ws.Range("E21").Value = 20.92914372033528
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 19.52106029856116
$ws.Range("H21").Value = 12.04383642256364
$ws.Range("I21").Value = 16.85225549085647
$ws.Range("K21").Value = 10.63077383990961
$ws.Range("O21").Value = 17.08277493906566

$ws.Range("B22").Value = 10.8339213143666
$ws.Range("C22").Value = 7.379666528304265
$ws.Range("E22").Value = 21.67573630976689
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 19.41801857517667
$ws.Range("H22").Value = 12.00089025307341
$ws.Range("I22").Value = 16.76480552755886
$ws.Range("K22").Value = 10.86816619793592
$ws.Range("O22").Value = 17.00524296824771

$ws.Range("B23").Value = 10.65964682820095
$ws.Range("C23").Value = 7.307482738245826
$ws.Range("E23").Value = 21.28049022332542
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 19.47220501867134
$ws.Range("H23").Value = 12.0236385960696
$ws.Range("I23").Value = 16.81117114908265
$ws.Range("K23").Value = 10.74222638918788
$ws.Range("O23").Value = 17.04625628587731

$ws.Range("B24").Value = 9.969284046364155
$ws.Range("C24").Value = 7.026203862639331
$ws.Range("E24").Value = 19.71253894986639
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 19.69507231651459
$ws.Range("H24").Value = 12.11359258207631
$ws.Range("I24").Value = 16.99354834499328
$ws.Range("K24").Value = 10.24876505866394
$ws.Range("O24").Value = 17.20963480470436

$ws.Range("B25").Value = 9.181356822046521
$ws.Range("C25").Value = 6.707046629126884
$ws.Range("E25").Value = 17.86835774785718
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 19.97237841847859
$ws.Range("H25").Value = 12.21874589428214
$ws.Range("I25").Value = 17.20479036690595
$ws.Range("K25").Value = 9.68292532718576
$ws.Range("O25").Value = 17.40296829654003
